# Adding the weathermap forecast for RAAL
# Update forecasted GHI/DNI/DHI values on the "Daily" and "Hourly" sheets.

$wb = $excel.ActiveWorkbook

# ----- Daily sheet (row 2) -----
$daily = $wb.Worksheets.Item("Daily")
$daily.Range("G2").Value = 2574.75
$daily.Range("H2").Value = 5761.09
$daily.Range("I2").Value = 673.88
$daily.Range("J2").Value = 662.85
$daily.Range("L2").Value = 662.85

# ----- Hourly sheet -----
$hourly = $wb.Worksheets.Item("Hourly")

# Row 9 (hour = 7)
$hourly.Range("I9").Value = 12.18
$hourly.Range("K9").Value = 0.5
$hourly.Range("M9").Value = 0.5

# Row 10 (hour = 8)
$hourly.Range("H10").Value = 76.05
$hourly.Range("I10").Value = 353.26
$hourly.Range("K10").Value = 23.83
$hourly.Range("M10").Value = 23.83

# Row 11 (hour = 9)
$hourly.Range("H11").Value = 211.67
$hourly.Range("I11").Value = 596.55
$hourly.Range("K11").Value = 52.92
$hourly.Range("M11").Value = 52.92

# Row 12 (hour = 10)
$hourly.Range("H12").Value = 330.44
$hourly.Range("I12").Value = 706.35
$hourly.Range("J12").Value = 82.06999999999999
$hourly.Range("K12").Value = 82.61
$hourly.Range("M12").Value = 82.61

# Row 13 (hour = 11)
$hourly.Range("H13").Value = 410.43
$hourly.Range("I13").Value = 759.5599999999999
$hourly.Range("J13").Value = 90
$hourly.Range("K13").Value = 102.61
$hourly.Range("M13").Value = 102.61

# Row 14 (hour = 12)
$hourly.Range("H14").Value = 441.25
$hourly.Range("I14").Value = 777.25
$hourly.Range("J14").Value = 92.77

# Row 15 (hour = 13)
$hourly.Range("H15").Value = 419.29
$hourly.Range("I15").Value = 764.83
$hourly.Range("J15").Value = 90.8
$hourly.Range("K15").Value = 104.82
$hourly.Range("M15").Value = 104.82

# Row 16 (hour = 14)
$hourly.Range("H16").Value = 347.11
$hourly.Range("I16").Value = 718.5599999999999
$hourly.Range("J16").Value = 83.8
$hourly.Range("K16").Value = 88.89
$hourly.Range("M16").Value = 88.89

# Row 17 (hour = 15)
$hourly.Range("H17").Value = 233.83
$hourly.Range("I17").Value = 621.21
$hourly.Range("K17").Value = 65.06
$hourly.Range("M17").Value = 65.06

# Row 18 (hour = 16)
$hourly.Range("H18").Value = 98.45999999999999
$hourly.Range("I18").Value = 412.56
$hourly.Range("K18").Value = 29.78
$hourly.Range("M18").Value = 29.78

# Row 19 (hour = 17)
$hourly.Range("I19").Value = 38.78
$hourly.Range("K19").Value = 1.51
$hourly.Range("M19").Value = 1.51
